# modifying variable names in data engineering
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "c1_1_8_deg_final"
$ws.Range("F1").Value = "c1_final_date"
$ws.Range("H1").Value = "c1_initial_date"
$ws.Range("I1").Value = "c2_1_5_deg_final"
$ws.Range("J1").Value = "c2_1_8_deg_final"
$ws.Range("K1").Value = "c2_2_deg_final"
$ws.Range("L1").Value = "c2_final_date"
$ws.Range("Q1").Value = "c2_initial_date"
